$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (e.g. "1.00")
# must be pre-formatted as Text so Excel keeps them as strings instead of
# coercing to a number (which would also drop formatting like trailing zeros).
$textCells = @(
    'D4',
    'D5',
    'D6',
    'D13',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D25',
    'D26',
    'D29',
    'D31',
    'D33',
    'D34',
    'D36',
    'D38',
    'D40',
    'D41',
    'D44',
    'D46',
    'D47',
    'D49',
    'D50',
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.798.54'
$ws.Range('E2').Value = '  +2.81%  '
$ws.Range('D3').Value = '2.436.56'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '566.29'
$ws.Range('E5').Value = '  +2.40%  '
$ws.Range('D6').Value = '167.18'
$ws.Range('E6').Value = '  +5.75%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('E9').Value = '  +9.59%  '
$ws.Range('D10').Value = '2.436.60'
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('E12').Value = '  +2.59%  '
$ws.Range('D13').Value = '4.72'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('E14').Value = '  +6.54%  '
$ws.Range('D15').Value = '69.705.47'
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('D16').Value = '2.883.60'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('E17').Value = '  +5.80%  '
$ws.Range('D18').Value = '2.440.77'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').Value = '10.89'
$ws.Range('E19').Value = '  +6.05%  '
$ws.Range('D20').Value = '344.55'
$ws.Range('E20').Value = '  +4.58%  '
$ws.Range('D21').Value = '7.20'
$ws.Range('E21').Value = '  +6.07%  '
$ws.Range('D22').Value = '3.87'
$ws.Range('E22').Value = '  +3.41%  '
$ws.Range('D23').Value = '2.01'
$ws.Range('E23').Value = '  +8.04%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '66.10'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('D26').Value = '3.84'
$ws.Range('E26').Value = '  +6.15%  '
$ws.Range('E27').Value = '  +5.94%  '
$ws.Range('D28').Value = '2.561.14'
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').Value = '0.0₃0856'
$ws.Range('E30').Value = '  +8.08%  '
$ws.Range('D31').Value = '7.39'
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('E32').Value = '  +10.75%  '
$ws.Range('D33').Value = '457.27'
$ws.Range('E33').Value = '  +9.20%  '
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').Value = '159.10'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('E37').Value = '  +7.88%  '
$ws.Range('D38').Value = '19.14'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '18.28'
$ws.Range('E40').Value = '  +3.91%  '
$ws.Range('D41').Value = '0.305'
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('E42').Value = '  +4.95%  '
$ws.Range('E43').Value = '  +5.10%  '
$ws.Range('D44').Value = '38.11'
$ws.Range('E44').Value = '  +2.34%  '
$ws.Range('E45').Value = '  +3.51%  '
$ws.Range('D46').Value = '2.10'
$ws.Range('E46').Value = '  +9.43%  '
$ws.Range('D47').Value = '136.18'
$ws.Range('E47').Value = '  +6.03%  '
$ws.Range('E48').Value = '  +3.97%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0723'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = '0.492'
$ws.Range('E50').Value = '  +3.94%  '
$ws.Range('E51').Value = '  +2.41%  '

# Restore default (Normal) style on the text-forced cells so no stray
# number-format styling is left behind on cells that should look unchanged.
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
